$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row (row 1): keep look&feel, just normalise the row height
#    (it no longer needs the extra height now that the rest of the sheet
#    follows the same 12.8 pt row height).
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 12.8

# ---------------------------------------------------------------------------
# 2. Existing data rows (2-5): renumber the id columns, keep dates/status,
#    and switch the font used in the date/status columns from "Times New
#    Roman" to "Arial" (this mirrors the style remap seen in the new file).
# ---------------------------------------------------------------------------
function Set-DataRow($r, $colaborador, $treinamento, $dataConclusao, $validade, $status) {
    $ws.Cells.Item($r, 1).Value = $colaborador
    $ws.Cells.Item($r, 2).Value = $treinamento

    $ws.Cells.Item($r, 3).Value = $dataConclusao
    $ws.Cells.Item($r, 4).Value = $validade
    $ws.Cells.Item($r, 5).Value = $status

    # columns A & B keep their existing (Times New Roman) font, centred, no wrap
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4108

    # columns C, D (dates) & E (status) move to Arial 10, centred
    foreach ($c in 3..5) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Font.Name = "Arial"
        $cell.Font.Size = 10
        $cell.HorizontalAlignment = -4108
    }
    $ws.Cells.Item($r, 5).WrapText = $true
}

Set-DataRow 2 11 11 44941 45672 "Completo"
Set-DataRow 3 12 12 44602 44967 "Vencido"
Set-DataRow 4 13 13 44275 45371 "Completo"
Set-DataRow 5 14 14 44691 45117 "Pendente"
Set-DataRow 6 15 15 44336 45087 "Vencido"
Set-DataRow 7 36 16 44941 45672 "Completo"

# ---------------------------------------------------------------------------
# 3. Column widths (cosmetic) matching the target sheet.
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 14.45
$ws.Columns(2).ColumnWidth = 13.83
$ws.Columns(3).ColumnWidth = 14.45

# ---------------------------------------------------------------------------
# 4. Selection / active cell follows the last populated cell (E7).
# ---------------------------------------------------------------------------
$ws.Range("E7").Select()

Write-Output "done"
